$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Insert a new row above row 8. This shifts the old row 8 (Desc) and everything
# below it down by one row, and inherits formatting from row 7 above.
$ws.Rows("8:8").Insert()

# Set the label for the newly inserted row 8 cell A8.
$ws.Range("A8").Value = "Force"

# Give the new row the same style as row 7 (copy formatting explicitly to be safe).
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the boolean formulas across B8:I8 just like the other similar rows.
$ws.Range("B8:I8").Formula = "=FALSE()"

$wb.Save()
